$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new data rows (50 and 51) ---
# Set string-valued cells first, in the same left-to-right / row-by-row
# order used originally so that new shared-string entries get appended
# to xl/sharedStrings.xml in the expected order:
#   R2DW, TN\L2_R2_R2XCENTOFS, TN\L_R2_R2XCENTOFS, r2_center_ofs
$ws.Cells.Item(50, 3).Value = "R2DW"
$ws.Cells.Item(51, 3).Value = "R2DW"
$ws.Cells.Item(50, 4).Value = "TN\L2_R2_R2XCENTOFS"
$ws.Cells.Item(51, 4).Value = "TN\L_R2_R2XCENTOFS"
$ws.Cells.Item(50, 2).Value = "r2_center_ofs"
$ws.Cells.Item(51, 2).Value = "r2_center_ofs"
$ws.Cells.Item(50, 1).Value = 1580
$ws.Cells.Item(51, 1).Value = 2250

# Match the fill/style used by the other rows in this block (same as row 2)
$ws.Range("A2:D2").Copy()
$ws.Range("A50:D51").PasteSpecial(-4122)

# --- Update the view state: scroll down and select D46 ---
$ws.Range("D46").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1

# --- Update the workbook window position ---
$excel.ActiveWindow.Left = 11160
